$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District (column G) names to official names
$ws.Range("G3").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G5").Value = "Davangere"
$ws.Range("G9").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G10").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G11").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G17").Value = "Davangere"
$ws.Range("G18").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G19").Value = "Bagalkot"
$ws.Range("G26").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G35").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G36").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G38").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G40").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G48").Value = "Chikkamagaluru (Chikmagalur)"

# Remove stray empty Address (column F) cells that had no content
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("F44").ClearContents()
$ws.Range("F50").ClearContents()
